$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 377.14285
$ws.Range("I19").Value = 446.66666
$ws.Range("J19").Value = 325
$ws.Range("K19").Value = 446.66666
$ws.Range("L19").Value = 325
$ws.Range("M19").Value = -271.66666
$ws.Range("N19").Value = -675

$ws.Range("H69").Value = 4653.2383
$ws.Range("I69").Value = 4000
$ws.Range("J69").Value = 5055.231
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 15165.693
$ws.Range("M69").Value = -11126
$ws.Range("N69").Value = -16913.693

$ws.Range("H72").Value = 4653.2383
$ws.Range("I72").Value = 4000
$ws.Range("J72").Value = 5055.231
$ws.Range("K72").Value = 36000
$ws.Range("L72").Value = 45497.079
$ws.Range("M72").Value = -31632
$ws.Range("N72").Value = -54233.079

$ws.Range("H80").Value = 484.5
$ws.Range("J80").Value = 660.6
$ws.Range("L80").Value = 1981.8
$ws.Range("N80").Value = -3977.8

$ws.Range("H83").Value = 484.5
$ws.Range("J83").Value = 660.6
$ws.Range("L83").Value = 5945.400000000001
$ws.Range("N83").Value = -15929.4

$ws.Range("H132").Value = 2239.0312
$ws.Range("I132").Value = 1703.1428
$ws.Range("J132").Value = 5990.25
$ws.Range("K132").Value = 5109.428400000001
$ws.Range("L132").Value = 17970.75
$ws.Range("M132").Value = -2579.428400000001
$ws.Range("N132").Value = -23030.75

$ws.Range("H138").Value = 2971.2656
$ws.Range("I138").Value = 1186.1818
$ws.Range("J138").Value = 4871.516
$ws.Range("K138").Value = 3558.5454
$ws.Range("L138").Value = 14614.548
$ws.Range("M138").Value = 1581.4546
$ws.Range("N138").Value = -24894.548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 146.66667
$ws.Range("I5").Value = 120
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 120
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = -8
$ws.Range("N5").Value = -424

$ws.Range("H61").Value = 3021.3809
$ws.Range("I61").Value = 2722.125
$ws.Range("J61").Value = 3979
$ws.Range("K61").Value = 2722.125
$ws.Range("L61").Value = 3979
$ws.Range("M61").Value = -2510.125
$ws.Range("N61").Value = -4403

$ws.Range("H74").Value = 1650.909
$ws.Range("I74").Value = 1477.9412
$ws.Range("J74").Value = 1834.6875
$ws.Range("K74").Value = 1477.9412
$ws.Range("L74").Value = 1834.6875
$ws.Range("M74").Value = -603.9412
$ws.Range("N74").Value = -3582.6875

$ws.Range("H77").Value = 1650.909
$ws.Range("I77").Value = 1477.9412
$ws.Range("J77").Value = 1834.6875
$ws.Range("K77").Value = 7389.706
$ws.Range("L77").Value = 9173.4375
$ws.Range("M77").Value = -3021.706
$ws.Range("N77").Value = -17909.4375

$ws.Range("H102").Value = 12346412
$ws.Range("I102").Value = 12346412
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 12346412
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -12344790
$ws.Range("N102").ClearContents()

$ws.Range("H132").Value = 2649.1345
$ws.Range("I132").Value = 2260.394
$ws.Range("J132").Value = 3324.3157
$ws.Range("K132").Value = 6781.181999999999
$ws.Range("L132").Value = 9972.947100000001
$ws.Range("M132").Value = -4251.181999999999
$ws.Range("N132").Value = -15032.9471

$ws.Range("H136").Value = 3021.3809
$ws.Range("I136").Value = 2722.125
$ws.Range("J136").Value = 3979
$ws.Range("K136").Value = 8166.375
$ws.Range("L136").Value = 11937
$ws.Range("M136").Value = -5616.375
$ws.Range("N136").Value = -17037

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 146.66667
$ws.Range("I4").Value = 120
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 120
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -5
$ws.Range("N4").Value = -430

$ws.Range("H134").Value = 3020.625
$ws.Range("I134").Value = 3041.25
$ws.Range("K134").Value = 9123.75
$ws.Range("M134").Value = -6588.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4429.3037
$ws.Range("I31").Value = 1543.8214
$ws.Range("J31").Value = 7314.7856
$ws.Range("K31").Value = 1543.8214
$ws.Range("L31").Value = 7314.7856
$ws.Range("M31").Value = -1248.8214
$ws.Range("N31").Value = -7904.7856

$ws.Range("H34").Value = 4429.3037
$ws.Range("I34").Value = 1543.8214
$ws.Range("J34").Value = 7314.7856
$ws.Range("K34").Value = 1543.8214
$ws.Range("L34").Value = 7314.7856
$ws.Range("M34").Value = -1341.8214
$ws.Range("N34").Value = -7718.7856

$ws.Range("H132").Value = 2038.3667
$ws.Range("I132").Value = 1746.4073
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 5239.2219
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -2709.2219
$ws.Range("N132").Value = -19058

$ws.Range("H134").Value = 2828.7354
$ws.Range("I134").Value = 2722.1
$ws.Range("K134").Value = 8166.299999999999
$ws.Range("M134").Value = -5631.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 14285852
$ws.Range("I23").Value = 33333346
$ws.Range("J23").Value = 231
$ws.Range("K23").Value = 100000038
$ws.Range("L23").Value = 693
$ws.Range("M23").Value = -99999803
$ws.Range("N23").Value = -1163

$ws.Range("H122").Value = 409.5263
$ws.Range("I122").Value = 352.53845
$ws.Range("J122").Value = 533
$ws.Range("K122").Value = 3172.84605
$ws.Range("L122").Value = 4797
$ws.Range("M122").Value = -722.8460500000001
$ws.Range("N122").Value = -9697

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4231
$ws.Range("I132").Value = 5210.4443
$ws.Range("J132").Value = 3847.739
$ws.Range("K132").Value = 15631.3329
$ws.Range("L132").Value = 11543.217
$ws.Range("M132").Value = -13101.3329
$ws.Range("N132").Value = -16603.217

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 7661.913
$ws.Range("I82").Value = 6652.6665
$ws.Range("J82").Value = 11295.2
$ws.Range("K82").Value = 6652.6665
$ws.Range("L82").Value = 11295.2
$ws.Range("M82").Value = -6291.6665
$ws.Range("N82").Value = -12017.2

$ws.Range("H85").Value = 7661.913
$ws.Range("I85").Value = 6652.6665
$ws.Range("J85").Value = 11295.2
$ws.Range("K85").Value = 6652.6665
$ws.Range("L85").Value = 11295.2
$ws.Range("M85").Value = -5404.6665
$ws.Range("N85").Value = -13791.2

$ws.Range("H93").Value = 41686016
$ws.Range("I93").Value = 40639.8
$ws.Range("J93").Value = 71432710
$ws.Range("K93").Value = 40639.8
$ws.Range("L93").Value = 71432710
$ws.Range("M93").Value = -39391.8
$ws.Range("N93").Value = -71435206

$ws.Range("H122").Value = 3703139.8
$ws.Range("I122").Value = 4203492
$ws.Range("J122").Value = 2001941
$ws.Range("K122").Value = 12610476
$ws.Range("L122").Value = 6005823
$ws.Range("M122").Value = -12608026
$ws.Range("N122").Value = -6010723

$ws.Range("H132").Value = 16053569
$ws.Range("I132").Value = 19700806
$ws.Range("J132").Value = 5732.2
$ws.Range("K132").Value = 59102418
$ws.Range("L132").Value = 17196.6
$ws.Range("M132").Value = -59099888
$ws.Range("N132").Value = -22256.6

$ws.Range("H136").Value = 5943.7
$ws.Range("I136").Value = 3622
$ws.Range("J136").Value = 19100
$ws.Range("K136").Value = 10866
$ws.Range("L136").Value = 57300
$ws.Range("M136").Value = -8316
$ws.Range("N136").Value = -62400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8591.916999999999
$ws.Range("I62").Value = 3100
$ws.Range("J62").Value = 14083.833
$ws.Range("K62").Value = 3100
$ws.Range("L62").Value = 14083.833
$ws.Range("M62").Value = -2476
$ws.Range("N62").Value = -15331.833

$ws.Range("H65").Value = 8591.916999999999
$ws.Range("I65").Value = 3100
$ws.Range("J65").Value = 14083.833
$ws.Range("K65").Value = 15500
$ws.Range("L65").Value = 70419.16500000001
$ws.Range("M65").Value = -12380
$ws.Range("N65").Value = -76659.16500000001

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H113").Value = 1155.6333
$ws.Range("I113").Value = 833.2857
$ws.Range("K113").Value = 2499.8571
$ws.Range("M113").Value = -329.8571000000002

$ws.Range("H122").Value = 2074.889
$ws.Range("I122").Value = 1856
$ws.Range("J122").Value = 2348.5
$ws.Range("K122").Value = 5568
$ws.Range("L122").Value = 7045.5
$ws.Range("M122").Value = -3118
$ws.Range("N122").Value = -11945.5

$ws.Range("H132").Value = 3366.7036
$ws.Range("I132").Value = 3117.6875
$ws.Range("J132").Value = 3728.9092
$ws.Range("K132").Value = 9353.0625
$ws.Range("L132").Value = 11186.7276
$ws.Range("M132").Value = -6823.0625
$ws.Range("N132").Value = -16246.7276

$ws.Range("H136").Value = 1335.1025
$ws.Range("I136").Value = 753.0833
$ws.Range("J136").Value = 2266.3333
$ws.Range("K136").Value = 2259.2499
$ws.Range("L136").Value = 6798.999899999999
$ws.Range("M136").Value = 290.7501000000002
$ws.Range("N136").Value = -11898.9999
